# deposito parte1 e 2 19/08/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transações")

# Linha 20 - deposito parte 1
$ws.Cells.Item(20, 1).Value = 567
$ws.Cells.Item(20, 2).Value = "Clevison"
$ws.Cells.Item(20, 3).Value = "678.234.987-89"
$ws.Cells.Item(20, 4).Value = "Depósito"
$ws.Cells.Item(20, 5).Value = 1000
$ws.Cells.Item(20, 6).Value = "19/08/2025"

# Linha 21 - deposito parte 2 (saque)
$ws.Cells.Item(21, 1).Value = 567
$ws.Cells.Item(21, 2).Value = "Clevison"
$ws.Cells.Item(21, 3).Value = "678.234.987-89"
$ws.Cells.Item(21, 4).Value = "Saque"
$ws.Cells.Item(21, 5).Value = 700
$ws.Cells.Item(21, 6).Value = "19/08/2025"
